# Update cryptos price/volume data (Sat Apr 15 17:43:31 UTC 2023 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.460.29'
$ws.Range('E2').Value = '  +0.74%  '
$ws.Range('D3').Value = '2.108.71'
$ws.Range('E3').Value = '  +1.67%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.007'
$ws.Range('E4').Value = '  +0.91%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '335.35'
$ws.Range('E5').Value = '  +2.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.007'
$ws.Range('E6').Value = '  +0.85%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5231'
$ws.Range('E7').Value = '  +0.86%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4549'
$ws.Range('E8').Value = '  +5.30%  '
$ws.Range('E9').Value = '  +15.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08938'
$ws.Range('E10').Value = '  +0.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.176'
$ws.Range('E11').Value = '  +2.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.33'
$ws.Range('E12').Value = '  +0.71%  '
$ws.Range('D13').Value = '2.103.13'
$ws.Range('E13').Value = '  +1.42%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.860'
$ws.Range('E14').Value = '  +3.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.054'
$ws.Range('E15').Value = '  +5.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '96.60'
$ws.Range('E16').Value = '  +1.91%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001144'
$ws.Range('E17').Value = '  +2.09%  '
$ws.Range('B18').Value = 'BinanceUSD'
$ws.Range('C18').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.008'
$ws.Range('E18').Value = '  +0.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06672'
$ws.Range('E19').Value = '  +1.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.23'
$ws.Range('E20').Value = '  +2.77%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.006'
$ws.Range('E21').Value = '  +0.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.361'
$ws.Range('E22').Value = '  +2.52%  '
$ws.Range('D23').Value = '30.511.32'
$ws.Range('E23').Value = '  +0.83%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.42'
$ws.Range('E24').Value = '  +1.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.373'
$ws.Range('E25').Value = '  +4.23%  '
$ws.Range('D26').Value = '2.353.55'
$ws.Range('E26').Value = '  +1.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.38'
$ws.Range('E27').Value = '  +0.99%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '163.75'
$ws.Range('E28').Value = '  +1.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.546'
$ws.Range('E29').Value = '  +0.82%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.57'
$ws.Range('E30').Value = '  +2.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.225'
$ws.Range('E31').Value = '  +3.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1073'
$ws.Range('E32').Value = '  +0.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.670'
$ws.Range('E33').Value = '  +3.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.331'
$ws.Range('E34').Value = '  +4.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.952'
$ws.Range('E35').Value = '  +3.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.48'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02585'
$ws.Range('E37').Value = '  +1.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.695'
$ws.Range('E38').Value = '  +5.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06844'
$ws.Range('E39').Value = '  +3.56%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2304'
$ws.Range('E40').Value = '  +3.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '12.66'
$ws.Range('E41').Value = '  +0.86%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6889'
$ws.Range('E42').Value = '  +1.81%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.249'
$ws.Range('E43').Value = '  +0.72%  '
$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.339'
$ws.Range('E44').Value = '  +6.76%  '
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.007'
$ws.Range('E45').Value = '  +0.86%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '14.02'
$ws.Range('E46').Value = '  +0.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6384'
$ws.Range('E47').Value = '  +0.78%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.674'
$ws.Range('E48').Value = '  +2.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.250'
$ws.Range('E49').Value = '  +1.44%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.3430'
$ws.Range('E50').Value = '  +26.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '83.38'
$ws.Range('E51').Value = '  +2.75%  '
